$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.047621086334526
$ws.Range("D2").Value = 1.052021915148968
$ws.Range("E2").Value = 1.060395046940329
$ws.Range("F2").Value = 1.066004923385548
$ws.Range("I2").Value = 1.040352457862874
$ws.Range("J2").Value = 1.052669180581595
$ws.Range("K2").Value = 1.054771799637858
$ws.Range("L2").Value = 1.06312194377474
$ws.Range("M2").Value = 1.068716637445639
$ws.Range("N2").Value = 1.054164091560863

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.048713905276257
$ws.Range("D3").Value = 1.05286247078866
$ws.Range("E3").Value = 1.061440337667215
$ws.Range("F3").Value = 1.067050173802101
$ws.Range("I3").Value = 1.040576026042302
$ws.Range("J3").Value = 1.053409842141715
$ws.Range("K3").Value = 1.055425046595825
$ws.Range("L3").Value = 1.063981068508169
$ws.Range("M3").Value = 1.069576822861148
$ws.Range("N3").Value = 1.054905804945354

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.049421146560831
$ws.Range("D4").Value = 1.053406098067983
$ws.Range("E4").Value = 1.062117118448942
$ws.Range("F4").Value = 1.067726818904605
$ws.Range("I4").Value = 1.040718845199083
$ws.Range("J4").Value = 1.053888639362997
$ws.Range("K4").Value = 1.055846827075686
$ws.Range("L4").Value = 1.06453678165307
$ws.Range("M4").Value = 1.070133111580824
$ws.Range("N4").Value = 1.055385282113609

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.049718498409071
$ws.Range("D5").Value = 1.053634574140532
$ws.Range("E5").Value = 1.062401734420752
$ws.Range("F5").Value = 1.068011351536451
$ws.Range("I5").Value = 1.040778444477682
$ws.Range("J5").Value = 1.054089815299003
$ws.Range("K5").Value = 1.056023924397159
$ws.Range("L5").Value = 1.064770355794876
$ws.Range("M5").Value = 1.070366901114289
$ws.Range("N5").Value = 1.055586743742519

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.049768426718858
$ws.Range("D6").Value = 1.053672932459673
$ws.Range("E6").Value = 1.062449528393769
$ws.Range("F6").Value = 1.068059129972207
$ws.Range("I6").Value = 1.040788425535896
$ws.Range("J6").Value = 1.054123587142864
$ws.Range("K6").Value = 1.056053646951946
$ws.Range("L6").Value = 1.064809571118357
$ws.Range("M6").Value = 1.070406151042429
$ws.Range("N6").Value = 1.055620563546273

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.04942511968199
$ws.Range("D7").Value = 1.053409151231943
$ws.Range("E7").Value = 1.062120921119418
$ws.Range("F7").Value = 1.067730620564913
$ws.Range("I7").Value = 1.040719643304343
$ws.Range("J7").Value = 1.053891327919833
$ws.Range("K7").Value = 1.055849194320534
$ws.Range("L7").Value = 1.06453990286981
$ws.Range("M7").Value = 1.070136235780162
$ws.Range("N7").Value = 1.055387974488505

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.047990386076797
$ws.Range("D8").Value = 1.052306040099171
$ws.Range("E8").Value = 1.06074822322586
$ws.Range("F8").Value = 1.066358108702876
$ws.Range("I8").Value = 1.040428395468697
$ws.Range("J8").Value = 1.052919585813474
$ws.Range("K8").Value = 1.05499275652369
$ws.Range("L8").Value = 1.063412329843366
$ws.Range("M8").Value = 1.069007404871271
$ws.Range("N8").Value = 1.054414852396893

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.045463066837111
$ws.Range("D9").Value = 1.050360186197967
$ws.Range("E9").Value = 1.058332489588247
$ws.Range("F9").Value = 1.063941866257704
$ws.Range("I9").Value = 1.039901068257001
$ws.Range("J9").Value = 1.051203735174068
$ws.Range("K9").Value = 1.053476619964723
$ws.Range("L9").Value = 1.061423896531759
$ws.Range("M9").Value = 1.067015909558175
$ws.Range("N9").Value = 1.052696565052769

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.043778749891438
$ws.Range("D10").Value = 1.049061615324242
$ws.Range("E10").Value = 1.056724122076982
$ws.Range("F10").Value = 1.062332600736602
$ws.Range("I10").Value = 1.039540047443261
$ws.Range("J10").Value = 1.050057479382208
$ws.Range("K10").Value = 1.052461186332669
$ws.Range("L10").Value = 1.060097273985101
$ws.Range("M10").Value = 1.065686681019354
$ws.Range("N10").Value = 1.051548681446207

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.043049550283973
$ws.Range("D11").Value = 1.048499008813488
$ws.Range("E11").Value = 1.056028186253334
$ws.Range("F11").Value = 1.06163614429564
$ws.Range("I11").Value = 1.039381478839742
$ws.Range("J11").Value = 1.049560580415565
$ws.Range("K11").Value = 1.052020386633369
$ws.Range("L11").Value = 1.05952259612209
$ws.Range("M11").Value = 1.065110741593259
$ws.Range("N11").Value = 1.051051076826036

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.042778710657393
$ws.Range("D12").Value = 1.04828998447464
$ws.Range("E12").Value = 1.055769759546719
$ws.Range("F12").Value = 1.061377504472053
$ws.Range("I12").Value = 1.03932224243523
$ws.Range("J12").Value = 1.049375925471467
$ws.Range("K12").Value = 1.051856487368658
$ws.Range("L12").Value = 1.059309098909409
$ws.Range("M12").Value = 1.064896755837646
$ws.Range("N12").Value = 1.050866159650737

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.042836805890452
$ws.Range("D13").Value = 1.048334823036478
$ws.Range("E13").Value = 1.055825189565085
$ws.Range("F13").Value = 1.061432981096881
$ws.Range("I13").Value = 1.039334964090246
$ws.Range("J13").Value = 1.049415538409805
$ws.Range("K13").Value = 1.051891651854969
$ws.Range("L13").Value = 1.059354896427085
$ws.Range("M13").Value = 1.064942659050874
$ws.Range("N13").Value = 1.050905828843991

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.043027162240692
$ws.Range("D14").Value = 1.048481731757876
$ws.Range("E14").Value = 1.05600682309413
$ws.Range("F14").Value = 1.06161476392239
$ws.Range("I14").Value = 1.039376589212891
$ws.Range("J14").Value = 1.049545318505025
$ws.Range("K14").Value = 1.052006842067247
$ws.Range("L14").Value = 1.059504949115601
$ws.Range("M14").Value = 1.065093054610116
$ws.Range("N14").Value = 1.051035793241833

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.043144449382601
$ws.Range("D15").Value = 1.048572240821807
$ws.Range("E15").Value = 1.056118743424775
$ws.Range("F15").Value = 1.061726773626044
$ws.Range("I15").Value = 1.03940219117962
$ws.Range("J15").Value = 1.049625269101722
$ws.Range("K15").Value = 1.052077792483747
$ws.Range("L15").Value = 1.059597396724278
$ws.Range("M15").Value = 1.065185710832076
$ws.Range("N15").Value = 1.051115857377548

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.043827146963547
$ws.Range("D16").Value = 1.049098947022523
$ws.Range("E16").Value = 1.056770319539533
$ws.Range("F16").Value = 1.062378829979713
$ws.Range("I16").Value = 1.03955052385368
$ws.Range("J16").Value = 1.05009044507486
$ws.Range("K16").Value = 1.052490417380254
$ws.Range("L16").Value = 1.060135408375137
$ws.Range("M16").Value = 1.06572489633516
$ws.Range("N16").Value = 1.051581693953924

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.0442554170478
$ws.Range("D17").Value = 1.049429251357791
$ws.Range("E17").Value = 1.057179169357776
$ws.Range("F17").Value = 1.062787945802924
$ws.Range("I17").Value = 1.039642968326778
$ws.Range("J17").Value = 1.050382087029496
$ws.Range("K17").Value = 1.052748948956949
$ws.Range("L17").Value = 1.060472824218305
$ws.Range("M17").Value = 1.066063012900573
$ws.Range("N17").Value = 1.051873750073587

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.044505231430306
$ws.Range("D18").Value = 1.049621881702513
$ws.Range("E18").Value = 1.057417692579549
$ws.Range("F18").Value = 1.063026611462649
$ws.Range("I18").Value = 1.039696673076955
$ws.Range("J18").Value = 1.050552142561422
$ws.Range("K18").Value = 1.052899638918816
$ws.Range("L18").Value = 1.060669609864009
$ws.Range("M18").Value = 1.066260194414399
$ws.Range("N18").Value = 1.052044047103874

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.044590413678643
$ws.Range("D19").Value = 1.049687558473196
$ws.Range("E19").Value = 1.057499030981394
$ws.Range("F19").Value = 1.063107996273056
$ws.Range("I19").Value = 1.039714948272586
$ws.Range("J19").Value = 1.050610117895659
$ws.Range("K19").Value = 1.052951002124788
$ws.Range("L19").Value = 1.060736704729781
$ws.Range("M19").Value = 1.066327422079264
$ws.Range("N19").Value = 1.052102104769736

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.044209466515825
$ws.Range("D20").Value = 1.049393815989508
$ws.Range("E20").Value = 1.057135298672184
$ws.Range("F20").Value = 1.062744047892583
$ws.Range("I20").Value = 1.039633072306641
$ws.Range("J20").Value = 1.050350802222866
$ws.Range("K20").Value = 1.052721222045393
$ws.Range("L20").Value = 1.060436625102066
$ws.Range("M20").Value = 1.066026739932288
$ws.Range("N20").Value = 1.051842420838944

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.042971106559434
$ws.Range("D21").Value = 1.048438472101982
$ws.Range("E21").Value = 1.055953334480554
$ws.Range("F21").Value = 1.061561231889065
$ws.Range("I21").Value = 1.039364340952027
$ws.Range("J21").Value = 1.049507103835595
$ws.Range("K21").Value = 1.051972926019223
$ws.Range("L21").Value = 1.059460763339774
$ws.Range("M21").Value = 1.065048768417413
$ws.Range("N21").Value = 1.050997524303189

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.042192601909103
$ws.Range("D22").Value = 1.047837536393033
$ws.Range("E22").Value = 1.055210619726588
$ws.Range("F22").Value = 1.060817867407421
$ws.Range("I22").Value = 1.039193429405899
$ws.Range("J22").Value = 1.048976148022492
$ws.Range("K22").Value = 1.051501478620277
$ws.Range("L22").Value = 1.058846990498162
$ws.Range("M22").Value = 1.064433553741375
$ws.Range("N22").Value = 1.05046581447193

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.042605292531119
$ws.Range("D23").Value = 1.048156129550273
$ws.Range("E23").Value = 1.055604305811003
$ws.Range("F23").Value = 1.061211908725625
$ws.Range("I23").Value = 1.039284217571723
$ws.Range("J23").Value = 1.04925766411591
$ws.Range("K23").Value = 1.051751493125771
$ws.Range("L23").Value = 1.059172382976334
$ws.Range("M23").Value = 1.064759721476411
$ws.Range("N23").Value = 1.050747730350491

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.044230229546931
$ws.Range("D24").Value = 1.04940982780156
$ws.Range("E24").Value = 1.057155121798031
$ws.Range("F24").Value = 1.062763883359426
$ws.Range("I24").Value = 1.039637544561217
$ws.Range("J24").Value = 1.050364938647945
$ws.Range("K24").Value = 1.052733750974155
$ws.Range("L24").Value = 1.060452981996399
$ws.Range("M24").Value = 1.066043130236897
$ws.Range("N24").Value = 1.051856577339367

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.046116338901742
$ws.Range("D25").Value = 1.050863473586564
$ws.Range("E25").Value = 1.058956641666746
$ws.Range("F25").Value = 1.064566248636024
$ws.Range("I25").Value = 1.040039065340935
$ws.Range("J25").Value = 1.051647739500347
$ws.Range("K25").Value = 1.053869403180086
$ws.Range("L25").Value = 1.061938131378479
$ws.Range("M25").Value = 1.06753103580453
$ws.Range("N25").Value = 1.053141199916119
